# Horarios Línea 141 - actualización de datos (scrap 08:39:08)
# Para cada hoja: se agregan las nuevas filas scrapeadas y se reordena
# todo el bloque de datos (desde la fila 6) por la columna B (Hora_Llegada),
# tal como hace el scraper al regenerar el archivo. Luego se actualizan los
# textos de cabecera "Última actualización" y "Total filas".

$wb = $excel.ActiveWorkbook

$nuevaHoraScrap = "08:39:08"

function Agregar-Filas {
    param(
        $ws,
        [object[]]$filas
    )

    $lastRow = $ws.Range("A5").End(4).Row
    $startRow = $lastRow + 1

    for ($i = 0; $i -lt $filas.Count; $i++) {
        $fila = $filas[$i]
        $r = $startRow + $i
        $ws.Cells.Item($r, 1).Value2 = $fila[0]
        $ws.Cells.Item($r, 2).Value2 = $fila[1]
        $ws.Cells.Item($r, 3).Value2 = $fila[2]
        $ws.Cells.Item($r, 4).Value2 = $fila[3]
        $ws.Cells.Item($r, 5).Value2 = $fila[4]
    }

    $newLastRow = $startRow + $filas.Count - 1
    $totalDataRows = $newLastRow - 6 + 1

    $rango = $ws.Range("A6:E$newLastRow")
    $clave = $ws.Range("B6:B$newLastRow")
    $rango.Sort($clave, 1)

    $ws.Range("A2").Value2 = "Última actualización: $nuevaHoraScrap"
    $ws.Range("A3").Value2 = "Total filas: $totalDataRows"
}

# --- Hoja LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$filas1 = @(
    @("08:39:08", "09:00", "215B_EL PATO", 21, "LP1912"),
    @("08:39:08", "09:04", "23_HERNANDEZ", 25, "LP1912"),
    @("08:39:08", "09:15", "11_ETCHEVERRY", 36, "LP1912"),
    @("08:39:08", "09:17", "27_EL RETIRO", 38, "LP1912"),
    @("08:39:08", "09:32", "23_HERNANDEZ", 53, "LP1912"),
    @("08:39:08", "09:45", "14_ABASTO", 66, "LP1912"),
    @("08:39:08", "10:05", "14_ABASTO", 86, "LP1912"),
    @("08:39:08", "10:15", "10_OLMOS", 96, "LP1912"),
    @("08:39:08", "10:30", "11_ETCHEVERRY", 111, "LP1912"),
    @("08:39:08", "10:34", "10_OLMOS", 115, "LP1912"),
    @("08:39:08", "10:37", "16_P MOR-SANTA ANA", 118, "LP1912")
)
Agregar-Filas $ws1 $filas1

# --- Hoja LP1912-215 ---
# (se usa el operador unario "," para evitar que PowerShell aplane
#  un arreglo-de-arreglos de un solo elemento)
$ws2 = $wb.Worksheets.Item("LP1912-215")
$filas2 = ,@("08:39:08", "09:00", "215B_EL PATO", 21, "LP1912")
Agregar-Filas $ws2 $filas2

# --- Hoja 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$filas3 = ,@("08:39:08", "10:31", "215B_LP-P MOR-1 Y 57", 112, "L6173")
Agregar-Filas $ws3 $filas3
